$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 318.51852
$ws.Range("I33").Value = 317.2
$ws.Range("J33").Value = 335
$ws.Range("K33").Value = 317.2
$ws.Range("L33").Value = 335
$ws.Range("M33").Value = -88.19999999999999
$ws.Range("N33").Value = -793

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3151.7856
$ws.Range("I64").Value = 2840
$ws.Range("J64").Value = 3219.5652
$ws.Range("K64").Value = 2840
$ws.Range("L64").Value = 3219.5652
$ws.Range("M64").Value = -2592
$ws.Range("N64").Value = -3715.5652

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3151.7856
$ws.Range("I67").Value = 2840
$ws.Range("J67").Value = 3219.5652
$ws.Range("K67").Value = 2840
$ws.Range("L67").Value = 3219.5652
$ws.Range("M67").Value = -1982
$ws.Range("N67").Value = -4935.5652

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 55668440
$ws.Range("I106").Value = 143994.14
$ws.Range("J106").Value = 250004000
$ws.Range("K106").Value = 143994.14
$ws.Range("L106").Value = 250004000
$ws.Range("M106").Value = -143363.14
$ws.Range("N106").Value = -250005262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 12790
$ws.Range("J37").Value = 12790
$ws.Range("L37").Value = 12790
$ws.Range("N37").Value = -13336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 21550
$ws.Range("J44").Value = 21550
$ws.Range("L44").Value = 21550
$ws.Range("N44").Value = -22526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 22393.334
$ws.Range("J55").Value = 22393.334
$ws.Range("L55").Value = 22393.334
$ws.Range("N55").Value = -23023.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1614.9333
$ws.Range("I110").Value = 1665.8182
$ws.Range("J110").Value = 1475
$ws.Range("K110").Value = 1665.8182
$ws.Range("L110").Value = 1475
$ws.Range("M110").Value = 379.1818000000001
$ws.Range("N110").Value = -5565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2234.25
$ws.Range("I107").Value = 2045.6666
$ws.Range("J107").Value = 2800
$ws.Range("K107").Value = 2045.6666
$ws.Range("L107").Value = 2800
$ws.Range("M107").Value = -125.6666
$ws.Range("N107").Value = -6640

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 851.9666999999999
$ws.Range("I16").Value = 791.1579
$ws.Range("J16").Value = 957
$ws.Range("K16").Value = 791.1579
$ws.Range("L16").Value = 957
$ws.Range("M16").Value = -504.1579
$ws.Range("N16").Value = -1531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277423.5
$ws.Range("I31").Value = 3143
$ws.Range("J31").Value = 8776667
$ws.Range("K31").Value = 3143
$ws.Range("L31").Value = 8776667
$ws.Range("M31").Value = -2848
$ws.Range("N31").Value = -8777257

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4277423.5
$ws.Range("I34").Value = 3143
$ws.Range("J34").Value = 8776667
$ws.Range("K34").Value = 3143
$ws.Range("L34").Value = 8776667
$ws.Range("M34").Value = -2941
$ws.Range("N34").Value = -8777071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 41668756
$ws.Range("I62").Value = 1803.3334
$ws.Range("K62").Value = 1803.3334
$ws.Range("M62").Value = -1179.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 41668756
$ws.Range("I65").Value = 1803.3334
$ws.Range("K65").Value = 9016.666999999999
$ws.Range("M65").Value = -5896.666999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2285.0417
$ws.Range("J99").Value = 2454.7273
$ws.Range("L99").Value = 2454.7273
$ws.Range("N99").Value = -5450.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 687.8444
$ws.Range("I107").Value = 738
$ws.Range("J107").Value = 635.4091
$ws.Range("K107").Value = 738
$ws.Range("L107").Value = 635.4091
$ws.Range("M107").Value = 1182
$ws.Range("N107").Value = -4475.4091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 851.9666999999999
$ws.Range("I113").Value = 791.1579
$ws.Range("J113").Value = 957
$ws.Range("K113").Value = 791.1579
$ws.Range("L113").Value = 957
$ws.Range("M113").Value = 1378.8421
$ws.Range("N113").Value = -5297

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 990.65216
$ws.Range("I122").Value = 899.75
$ws.Range("J122").Value = 1198.4286
$ws.Range("K122").Value = 2699.25
$ws.Range("L122").Value = 3595.2858
$ws.Range("M122").Value = -249.25
$ws.Range("N122").Value = -8495.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2285.0417
$ws.Range("J126").Value = 2454.7273
$ws.Range("L126").Value = 7364.1819
$ws.Range("N126").Value = -12304.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3292030.8
$ws.Range("I132").Value = 2224.5908
$ws.Range("J132").Value = 7815514
$ws.Range("K132").Value = 6673.7724
$ws.Range("L132").Value = 23446542
$ws.Range("M132").Value = -4143.7724
$ws.Range("N132").Value = -23451602

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1168.6842
$ws.Range("I5").Value = 360.4
$ws.Range("J5").Value = 1457.3572
$ws.Range("K5").Value = 1081.2
$ws.Range("L5").Value = 4372.071599999999
$ws.Range("M5").Value = -969.1999999999998
$ws.Range("N5").Value = -4596.071599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 801
$ws.Range("I122").Value = 406.57144
$ws.Range("J122").Value = 1491.25
$ws.Range("K122").Value = 3659.14296
$ws.Range("L122").Value = 13421.25
$ws.Range("M122").Value = -1209.14296
$ws.Range("N122").Value = -18321.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 739.5700000000001
$ws.Range("I131").Value = 385.36365
$ws.Range("J131").Value = 783.3483
$ws.Range("K131").Value = 1156.09095
$ws.Range("L131").Value = 2350.0449
$ws.Range("M131").Value = 3883.90905
$ws.Range("N131").Value = -12430.0449

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1168.6842
$ws.Range("I135").Value = 360.4
$ws.Range("J135").Value = 1457.3572
$ws.Range("K135").Value = 3243.6
$ws.Range("L135").Value = 13116.2148
$ws.Range("M135").Value = -708.5999999999999
$ws.Range("N135").Value = -18186.2148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 9643.333000000001
$ws.Range("J48").Value = 9643.333000000001
$ws.Range("L48").Value = 9643.333000000001
$ws.Range("N48").Value = -10613.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 12000
$ws.Range("J49").Value = 12000
$ws.Range("L49").Value = 12000
$ws.Range("N49").Value = -12368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("N74").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("N77").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 307.31818
$ws.Range("I107").Value = 378.9
$ws.Range("J107").Value = 247.66667
$ws.Range("K107").Value = 378.9
$ws.Range("L107").Value = 247.66667
$ws.Range("M107").Value = 1541.1
$ws.Range("N107").Value = -4087.66667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 20835422
$ws.Range("I61").Value = 2137
$ws.Range("J61").Value = 66668650
$ws.Range("K61").Value = 2137
$ws.Range("L61").Value = 66668650
$ws.Range("M61").Value = -1935
$ws.Range("N61").Value = -66669054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 20835422
$ws.Range("I113").Value = 2137
$ws.Range("J113").Value = 66668650
$ws.Range("K113").Value = 2137
$ws.Range("L113").Value = 66668650
$ws.Range("M113").Value = 33
$ws.Range("N113").Value = -66672990

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4025.6667
$ws.Range("I122").Value = 4538.5
$ws.Range("K122").Value = 13615.5
$ws.Range("M122").Value = -11165.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2007.7693
$ws.Range("I132").Value = 1743.1428
$ws.Range("J132").Value = 2316.5
$ws.Range("K132").Value = 5229.428400000001
$ws.Range("L132").Value = 6949.5
$ws.Range("M132").Value = -2699.428400000001
$ws.Range("N132").Value = -12009.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 385.32
$ws.Range("I107").Value = 290
$ws.Range("J107").Value = 528.3
$ws.Range("K107").Value = 870
$ws.Range("L107").Value = 1584.9
$ws.Range("M107").Value = 1050
$ws.Range("N107").Value = -5424.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2395.6956
$ws.Range("I122").Value = 2137.4546
$ws.Range("J122").Value = 2632.4167
$ws.Range("K122").Value = 6412.3638
$ws.Range("L122").Value = 7897.250100000001
$ws.Range("M122").Value = -3962.3638
$ws.Range("N122").Value = -12797.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1346.4634
$ws.Range("I132").Value = 1211.9117
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3635.7351
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1105.7351
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2344.7407
$ws.Range("I136").Value = 2385.4
$ws.Range("J136").Value = 2228.5715
$ws.Range("K136").Value = 7156.200000000001
$ws.Range("L136").Value = 6685.7145
$ws.Range("M136").Value = -4606.200000000001
$ws.Range("N136").Value = -11785.7145
